# Complete Automation of NewTours Site
# Adds the "excelbookflight" worksheet (flight booking confirmation data)
# after the existing "excelflightfind" sheet, making it the active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the last existing sheet so it lands at
# the end of the tab strip (sheetId 3 / rId3), and make it the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "excelbookflight"

# Passenger / booking details typed in by the automated NewTours flow.
$ws.Range("A1").Value = "Imran"
$ws.Range("B1").Value = "Hassan"
$ws.Range("C1").Value = "Tasrifa"
$ws.Range("D1").Value = "Pomy"
$ws.Range("E1").Value = "BA"
$ws.Range("F1").Value = 123456789
$ws.Range("G1").Value = "Imran"
$ws.Range("H1").Value = "Hassan"
$ws.Range("I1").Value = "Aiub"
$ws.Range("J1").Value = "Adabor"
$ws.Range("K1").Value = "Dhaka"
$ws.Range("L1").Value = "Adabor"
$ws.Range("M1").Value = 1207
$ws.Range("N1").Value = 200

# Column F was widened (best-fit) to comfortably show the card number.
$ws.Columns.Item(6).ColumnWidth = 9.81640625

# Leave the cursor where the automation script left it.
$ws.Range("M14").Select() | Out-Null

# Match the printed page setup captured when the booking sheet was saved.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
